# Insert a new data row at row 275 (pushing the existing rows 275-369 down to
# 276-370). The newly inserted row receives a brand-new price record while all
# the other rows simply shift down by one, so the last existing row (old 369)
# ends up duplicated into the new last row (370).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 275..369 down by one row, creating space for the new record and
# growing the sheet to row 370. Excel copies formatting/styles along with the
# cells being pushed down, so row 370 ends up identical to the old row 369.
$ws.Rows.Item(275).Insert()

# Populate the newly inserted row 275 with the new weekly price record.
$ws.Cells.Item(275, 1).Value2 = 5
$ws.Cells.Item(275, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(275, 3).Value2 = "Maule"

# Column D (Fecha) keeps the date-time number format used by the other rows.
$ws.Cells.Item(275, 4).NumberFormat = $ws.Cells.Item(276, 4).NumberFormat
$ws.Cells.Item(275, 4).Value2 = 44795

$ws.Cells.Item(275, 5).Value2 = 7
$ws.Cells.Item(275, 6).Value2 = 100114014
$ws.Cells.Item(275, 7).Value2 = "Betarraga"
$ws.Cells.Item(275, 8).Value2 = "Sin especificar"
$ws.Cells.Item(275, 9).Value2 = "Primera"
$ws.Cells.Item(275, 10).Value2 = 4000
$ws.Cells.Item(275, 11).Value2 = 750
$ws.Cells.Item(275, 12).Value2 = 750
$ws.Cells.Item(275, 13).Value2 = 750
$ws.Cells.Item(275, 14).Value2 = "$/paquete 5 unidades"
$ws.Cells.Item(275, 15).Value2 = "Región del Maule"
$ws.Cells.Item(275, 16).Value2 = 150
$ws.Cells.Item(275, 17).Value2 = 5
$ws.Cells.Item(275, 18).Value2 = "Hortaliza"
